# Zeiten.xlsx - "Sourcecode ausdokumentiert und Zeitenliste angepasst"
#
# Adds five new time-tracking rows (10-14) to the sheet, reusing the
# date/time/number styles of the last existing row (row 9) so no new
# cellXfs entries are created, and widens column E to fit the longest
# new task description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (date style, time styles) of the last data row
# into the five new rows so the new cells pick up the same styles
# (s="1" for dates, s="3" for times) instead of minting new ones.
$ws.Range("A9:E9").Copy($ws.Range("A10:E14"))

# Row 10 - 24.03.2015, Tommy, UML architecture
$ws.Range("A10").Value = 42087
$ws.Range("B10").Value = 0.83333333333333337
$ws.Range("C10").Value = 0.91666666666666663
$ws.Range("D10").Value = "Tommy"

# Row 11 - 26.03.2015, Juliano/Tommy/Constantin, difficulty-selection fragment
$ws.Range("A11").Value = 42089
$ws.Range("B11").Value = 0.46875
$ws.Range("C11").Value = 0.625
$ws.Range("D11").Value = "Juliano, Tommy, Constantin"
$ws.Range("E11").Value = "Fragment zur Schwierigkeitsauswahl erstellt und die erste View zum zeichnen der Schlange implementiert"

# Row 12 - 27.03.2015, Constantin, speed passed to SpieleActivty
$ws.Range("A12").Value = 42090
$ws.Range("B12").Value = 0.75
$ws.Range("C12").Value = 0.79166666666666663
$ws.Range("D12").Value = "Constantin"
$ws.Range("E12").Value = "Ausgwählte Geschwindigkeit in die SpieleActivty übergeben"

# Back to row 10's task text (entered after rows 11/12 by the original author)
$ws.Range("E10").Value = "Erstellen der UML-Architektur"

# Row 13 - 29.03.2015, Tommy, snake segments
$ws.Range("A13").Value = 42092
$ws.Range("B13").Value = 0.91666666666666663
$ws.Range("C13").Value = 0.95833333333333337
$ws.Range("D13").Value = "Tommy"
$ws.Range("E13").Value = "Schlangenglieder eingebaut und Activity dynamisiert"

# Row 14 - 29.03.2015, Juliano, source code documented
$ws.Range("A14").Value = 42092
$ws.Range("B14").Value = 0.91666666666666663
$ws.Range("C14").Value = 0.95833333333333337
$ws.Range("D14").Value = "Juliano"
$ws.Range("E14").Value = "Sourcecode ausdokumentiert"

# Widen column E to fit the new (longer) task descriptions.
$ws.Columns.Item(5).ColumnWidth = 96.25

# Leave the selection where the author's cursor ended up (the empty row
# right after the new data).
[void]$ws.Range("D15").Select()
